$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at position 7 (pushes existing rows 7.. down by one,
#    carrying over formatting/formulas, matching Excel's normal "insert row above" UX).
$ws.Rows.Item(7).Insert()

# 2) Populate the newly inserted row 7 with the new logging-notification entry
#    "A 34759-2023" that the diff adds.
$ws.Range("A7").Value = "A 34759-2023"
$ws.Range("B7").Value = 45139
$ws.Range("C7").Value = 45189
$ws.Range("D7").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E7").Value = "SURAHAMMAR"
$ws.Range("F7").Value = "Bergvik skog väst AB"
$ws.Range("G7").Value = 38.3
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 5
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 13
$ws.Range("R7").Value = "Blå taggsvamp`r`nGrantaggsvamp`r`nMotaggsvamp`r`nSkogshare`r`nUllticka`r`nDropptaggsvamp`r`nGrönpyrola`r`nMindre märgborre`r`nPlattlummer`r`nVedticka`r`nLopplummer`r`nMattlummer`r`nRevlummer"

$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/artfynd/A 34759-2023.xlsx", "A 34759-2023")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/kartor/A 34759-2023.png", "A 34759-2023")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomål/A 34759-2023.docx", "A 34759-2023")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomålsmail/A 34759-2023.docx", "A 34759-2023")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsyn/A 34759-2023.docx", "A 34759-2023")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsynsmail/A 34759-2023.docx", "A 34759-2023")'

# 3) The previous entry for "A 34759-2023" (old row 10, now shifted to row 11 by the
#    insert above) has been superseded by the new row-7 data above, so remove that
#    now-duplicate row. This shifts rows 12.. back up by one, restoring the original
#    198-data-row / 199-total-row layout with every other notification unchanged.
$ws.Rows.Item(11).Delete()

# 4) The workbook's "Förändrad" (last-refreshed) column C is stamped with today's
#    export date for every notification row; bump it from 45188 to 45189 throughout.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3)).Value = 45189
